$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Remove the display_name column (column C) entirely; remaining columns shift left.
$ws.Range("C1").EntireColumn.Delete()

# Re-express the column freeze (was frozen at column D/row 2, i.e. 3 cols x 1 row)
# to now freeze at column C/row 2 (2 cols x 1 row), matching the removed column,
# and leave the former display_name column (now showing label::language) selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C2").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("C1:C1048576").Select() | Out-Null

# Restore original active sheet (the edit above switches the active tab to
# "choices"; the source workbook had "survey" selected/active).
$wb.Worksheets.Item("survey").Activate()
